{"js": "// Replace the 25 three-digit-by-one-digit multiplication equations in the\n// practice table with their new values, cell by cell, preserving all\n// existing run/paragraph formatting (font, size, justification, etc.).\nconst replacements = [\n  [\"382\u00d79=3438\", \"286\u00d72=572\"],\n  [\"452\u00d74=1808\", \"357\u00d77=2499\"],\n  [\"529\u00d76=3174\", \"103\u00d75=515\"],\n  [\"524\u00d79=4716\", \"734\u00d74=2936\"],\n  [\"947\u00d76=5682\", \"518\u00d77=3626\"],\n  [\"437\u00d74=1748\", \"551\u00d77=3857\"],\n  [\"909\u00d77=6363\", \"686\u00d78=5488\"],\n  [\"908\u00d77=6356\", \"969\u00d79=8721\"],\n  [\"736\u00d75=3680\", \"381\u00d79=3429\"],\n  [\"158\u00d74=632\", \"481\u00d75=2405\"],\n  [\"617\u00d72=1234\", \"717\u00d77=5019\"],\n  [\"993\u00d78=7944\", \"405\u00d72=810\"],\n  [\"152\u00d79=1368\", \"850\u00d73=2550\"],\n  [\"664\u00d75=3320\", \"587\u00d73=1761\"],\n  [\"551\u00d79=4959\", \"249\u00d75=1245\"],\n  [\"221\u00d72=442\", \"379\u00d75=1895\"],\n  [\"270\u00d72=540\", \"161\u00d76=966\"],\n  [\"329\u00d74=1316\", \"534\u00d73=1602\"],\n  [\"307\u00d74=1228\", \"130\u00d78=1040\"],\n  [\"683\u00d79=6147\", \"992\u00d72=1984\"],\n  [\"838\u00d76=5028\", \"881\u00d78=7048\"],\n  [\"606\u00d75=3030\", \"548\u00d76=3288\"],\n  [\"735\u00d78=5880\", \"674\u00d74=2696\"],\n  [\"611\u00d79=5499\", \"341\u00d73=1023\"],\n  [\"612\u00d79=5508\", \"514\u00d72=1028\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication equations in the\n# practice table with their new values, preserving all existing\n# run/paragraph formatting (font, size, justification, etc.).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"382\u00d79=3438\", \"286\u00d72=572\"),\n    @(\"452\u00d74=1808\", \"357\u00d77=2499\"),\n    @(\"529\u00d76=3174\", \"103\u00d75=515\"),\n    @(\"524\u00d79=4716\", \"734\u00d74=2936\"),\n    @(\"947\u00d76=5682\", \"518\u00d77=3626\"),\n    @(\"437\u00d74=1748\", \"551\u00d77=3857\"),\n    @(\"909\u00d77=6363\", \"686\u00d78=5488\"),\n    @(\"908\u00d77=6356\", \"969\u00d79=8721\"),\n    @(\"736\u00d75=3680\", \"381\u00d79=3429\"),\n    @(\"158\u00d74=632\",  \"481\u00d75=2405\"),\n    @(\"617\u00d72=1234\", \"717\u00d77=5019\"),\n    @(\"993\u00d78=7944\", \"405\u00d72=810\"),\n    @(\"152\u00d79=1368\", \"850\u00d73=2550\"),\n    @(\"664\u00d75=3320\", \"587\u00d73=1761\"),\n    @(\"551\u00d79=4959\", \"249\u00d75=1245\"),\n    @(\"221\u00d72=442\",  \"379\u00d75=1895\"),\n    @(\"270\u00d72=540\",  \"161\u00d76=966\"),\n    @(\"329\u00d74=1316\", \"534\u00d73=1602\"),\n    @(\"307\u00d74=1228\", \"130\u00d78=1040\"),\n    @(\"683\u00d79=6147\", \"992\u00d72=1984\"),\n    @(\"838\u00d76=5028\", \"881\u00d78=7048\"),\n    @(\"606\u00d75=3030\", \"548\u00d76=3288\"),\n    @(\"735\u00d78=5880\", \"674\u00d74=2696\"),\n    @(\"611\u00d79=5499\", \"341\u00d73=1023\"),\n    @(\"612\u00d79=5508\", \"514\u00d72=1028\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $old,    # FindText\n        $true,   # MatchCase\n        $false,  # MatchWholeWord\n        $false,  # MatchWildcards\n        $false,  # MatchSoundsLike\n        $false,  # MatchAllWordForms\n        $true,   # Forward\n        1,       # Wrap (wdFindContinue)\n        $false,  # Format\n        $new,    # ReplaceWith\n        2        # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n"}
